$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.082.88'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.592.34'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.55'
$ws.Range("E5").Value = '  -1.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '191.75'
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.587.75'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +2.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.665'
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("E12").Value = '  -3.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000306'
$ws.Range("E13").Value = '  +5.35%  '
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.170.60'
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.00'
$ws.Range("E16").Value = '  +3.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.591.76'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.096.71'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.69'
$ws.Range("E19").Value = '  +1.77%  '
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '477.75'
$ws.Range("E22").Value = '  -3.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.62'
$ws.Range("E23").Value = '  +10.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.05'
$ws.Range("E24").Value = '  -6.17%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.39'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '95.52'
$ws.Range("E26").Value = '  +5.33%  '
$ws.Range("E27").Value = '  -2.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.10'
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.38'
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.24'
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("E33").Value = '  +1.52%  '
$ws.Range("E34").Value = '  +2.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '586.98'
$ws.Range("E35").Value = '  -5.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.14'
$ws.Range("E36").Value = '  +2.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0808'
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.397'
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("E40").Value = '  +20.23%  '
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("E43").Value = '  +7.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.236.24'
$ws.Range("E44").Value = '  -2.40%  '
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.36'
$ws.Range("E47").Value = '  +2.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.45'
$ws.Range("E48").Value = '  +3.38%  '
$ws.Range("E49").Value = '  +0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("E51").Value = '  -4.88%  '
